# Felretting og tillegg i flyttall-listeeksempel
#
# Slide 1 has a "Tabell 5" table (3rd shape on the slide) that documents
# Python variable examples. Two cells need fixing:
#   1) The flyttall-liste "l" had a typo (comma instead of punctum):
#        l = [2.3, 2.5, 3.3, 1.9, 1,4]  ->  l = [2.3, 2.5, 3.3, 1.9, 1.4]
#   2) The explanation cell's last line documented indexing with l[2],
#      which no longer matches; it is corrected to use negative indexing
#      and the correct resulting value:
#        l[2] gir 2.5  ->  l[-1] gir 1.4

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(3)
$tbl = $sh.Table

# 1) Fix the typo in the list literal (row 7, column 1).
$listCell = $tbl.Cell(7, 1)
$listCell.Shape.TextFrame.TextRange.Text = "l = [2.3, 2.5, 3.3, 1.9, 1.4]"

# 2) Update the last explanatory line (row 8, column 1, last paragraph)
#    from "l[2] gir 2.5" to "l[-1] gir 1.4".
$explCell = $tbl.Cell(8, 1)
$explRange = $explCell.Shape.TextFrame.TextRange
$paragraphs = $explRange.Paragraphs()
$lastParagraph = $paragraphs.Item($paragraphs.Count)
$lastParagraph.Text = "l[-1] gir 1.4"
